$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.6648458909848785
$ws.Cells.Item(2, 3).Value = 0.1469525115199595
$ws.Cells.Item(2, 5).Value = 0.1125772746831366
$ws.Cells.Item(2, 6).Value = 0.4443680307746263
$ws.Cells.Item(2, 7).Value = 0.002470104542915695
$ws.Cells.Item(2, 11).Value = 0.3269720627495758
$ws.Cells.Item(2, 12).Value = 0.195439465471182
$ws.Cells.Item(2, 14).Value = 1.814608252324877
$ws.Cells.Item(2, 15).Value = 3.533118245477112

$ws.Cells.Item(3, 2).Value = 0.6261141630804445
$ws.Cells.Item(3, 3).Value = 0.1472037131370421
$ws.Cells.Item(3, 5).Value = 0.111892959911799
$ws.Cells.Item(3, 6).Value = 0.387822817061874
$ws.Cells.Item(3, 7).Value = 0.002472381189778544
$ws.Cells.Item(3, 11).Value = 0.2943811997274395
$ws.Cells.Item(3, 12).Value = 0.1884937971703522
$ws.Cells.Item(3, 14).Value = 1.833039586622277
$ws.Cells.Item(3, 15).Value = 3.558713352487473

$ws.Cells.Item(4, 2).Value = 0.6025688582830355
$ws.Cells.Item(4, 3).Value = 0.1473764394070756
$ws.Cells.Item(4, 5).Value = 0.1115303277697279
$ws.Cells.Item(4, 6).Value = 0.3531389305168915
$ws.Cells.Item(4, 7).Value = 0.002473854066225555
$ws.Cells.Item(4, 11).Value = 0.274411253764967
$ws.Cells.Item(4, 12).Value = 0.184331455755455
$ws.Cells.Item(4, 14).Value = 1.844950367903236
$ws.Cells.Item(4, 15).Value = 3.576197080501672

$ws.Cells.Item(5, 2).Value = 0.5930339779926896
$ws.Cells.Item(5, 3).Value = 0.1474514953062354
$ws.Cells.Item(5, 5).Value = 0.1113970457058393
$ws.Cells.Item(5, 6).Value = 0.3390132514313251
$ws.Cells.Item(5, 7).Value = 0.002474473191142673
$ws.Cells.Item(5, 11).Value = 0.2662841092043209
$ws.Cells.Item(5, 12).Value = 0.1826610834216069
$ws.Cells.Item(5, 14).Value = 1.84995348090969
$ws.Cells.Item(5, 15).Value = 3.583766527152676

$ws.Cells.Item(6, 2).Value = 0.5914543625962665
$ws.Cells.Item(6, 3).Value = 0.147464240831372
$ws.Cells.Item(6, 5).Value = 0.1113757904148329
$ws.Cells.Item(6, 6).Value = 0.336668177824194
$ws.Cells.Item(6, 7).Value = 0.002474577140401424
$ws.Cells.Item(6, 11).Value = 0.2649352673090988
$ws.Cells.Item(6, 12).Value = 0.1823852810577904
$ws.Cells.Item(6, 14).Value = 1.850793266991783
$ws.Cells.Item(6, 15).Value = 3.585050287431756

$ws.Cells.Item(7, 2).Value = 0.6024400236685779
$ws.Cells.Item(7, 3).Value = 0.1473774327071027
$ws.Cells.Item(7, 5).Value = 0.1115284715647995
$ws.Cells.Item(7, 6).Value = 0.3529483938344953
$ws.Cells.Item(7, 7).Value = 0.002473862339240373
$ws.Cells.Item(7, 11).Value = 0.2743016040230515
$ws.Cells.Item(7, 12).Value = 0.1843088238913566
$ws.Cells.Item(7, 14).Value = 1.845017236775988
$ws.Cells.Item(7, 15).Value = 3.576297364126731

$ws.Cells.Item(8, 2).Value = 0.6514425557570576
$ws.Cells.Item(8, 3).Value = 0.1470353005340606
$ws.Cells.Item(8, 5).Value = 0.112329394817273
$ws.Cells.Item(8, 6).Value = 0.4248636149813336
$ws.Cells.Item(8, 7).Value = 0.002470873996435944
$ws.Cells.Item(8, 11).Value = 0.3157264985513564
$ws.Cells.Item(8, 12).Value = 0.1930234093019152
$ws.Cells.Item(8, 14).Value = 1.8208401448441
$ws.Cells.Item(8, 15).Value = 3.54157647057653

$ws.Cells.Item(9, 2).Value = 0.749386863770269
$ws.Cells.Item(9, 3).Value = 0.1465102179588413
$ws.Cells.Item(9, 5).Value = 0.1143556909950334
$ws.Cells.Item(9, 6).Value = 0.5661985755041457
$ws.Cells.Item(9, 7).Value = 0.002465606492366814
$ws.Cells.Item(9, 11).Value = 0.3972692395071249
$ws.Cells.Item(9, 12).Value = 0.2109219611962345
$ws.Cells.Item(9, 14).Value = 1.778138552659412
$ws.Cells.Item(9, 15).Value = 3.487518035668586

$ws.Cells.Item(10, 2).Value = 0.8224512922350584
$ws.Cells.Item(10, 3).Value = 0.1462122212854169
$ws.Cells.Item(10, 5).Value = 0.1161213665064338
$ws.Cells.Item(10, 6).Value = 0.6702781546542269
$ws.Cells.Item(10, 7).Value = 0.002462094214007186
$ws.Cells.Item(10, 11).Value = 0.457350660992347
$ws.Cells.Item(10, 12).Value = 0.2245634210957377
$ws.Cells.Item(10, 14).Value = 1.749632459151938
$ws.Cells.Item(10, 15).Value = 3.456354362015276

$ws.Cells.Item(11, 2).Value = 0.8559255833357327
$ws.Cells.Item(11, 3).Value = 0.1460954805178858
$ws.Cells.Item(11, 5).Value = 0.1169845967248762
$ws.Cells.Item(11, 6).Value = 0.7176906081379002
$ws.Cells.Item(11, 7).Value = 0.002460573321988853
$ws.Cells.Item(11, 11).Value = 0.4847175171409503
$ws.Cells.Item(11, 12).Value = 0.2308756516770245
$ws.Cells.Item(11, 14).Value = 1.737285874855189
$ws.Cells.Item(11, 15).Value = 3.444034722061133

$ws.Cells.Item(12, 2).Value = 0.8686349420553086
$ws.Cells.Item(12, 3).Value = 0.1460539606656006
$ws.Cells.Item(12, 5).Value = 0.1173200896216677
$ws.Cells.Item(12, 6).Value = 0.7356546913071611
$ws.Cells.Item(12, 7).Value = 0.002460008395865355
$ws.Cells.Item(12, 11).Value = 0.4950853443176584
$ws.Cells.Item(12, 12).Value = 0.2332812030024485
$ws.Cells.Item(12, 14).Value = 1.732699808629636
$ws.Cells.Item(12, 15).Value = 3.43963662841125

$ws.Cells.Item(13, 2).Value = 0.865896281166556
$ws.Cells.Item(13, 3).Value = 0.1460627834876504
$ws.Cells.Item(13, 5).Value = 0.1172474528344942
$ws.Cells.Item(13, 6).Value = 0.7317853510981394
$ws.Cells.Item(13, 7).Value = 0.002460129574170812
$ws.Cells.Item(13, 11).Value = 0.4928522511627023
$ws.Cells.Item(13, 12).Value = 0.2327624478458574
$ws.Cells.Item(13, 14).Value = 1.733683527270879
$ws.Cells.Item(13, 15).Value = 3.440571956688615

$ws.Cells.Item(14, 2).Value = 0.8569705251134678
$ws.Cells.Item(14, 3).Value = 0.1460920109051216
$ws.Cells.Item(14, 5).Value = 0.1170120255801059
$ws.Cells.Item(14, 6).Value = 0.7191683204515869
$ws.Cells.Item(14, 7).Value = 0.002460526625131003
$ws.Cells.Item(14, 11).Value = 0.4855703954255546
$ws.Cells.Item(14, 12).Value = 0.2310732527922283
$ws.Cells.Item(14, 14).Value = 1.736906786244566
$ws.Cells.Item(14, 15).Value = 3.443667535028283

$ws.Cells.Item(15, 2).Value = 0.8515075670225656
$ws.Cells.Item(15, 3).Value = 0.1461102629554389
$ws.Cells.Item(15, 5).Value = 0.1168689396807245
$ws.Cells.Item(15, 6).Value = 0.7114413442032514
$ws.Cells.Item(15, 7).Value = 0.002460771261430521
$ws.Cells.Item(15, 11).Value = 0.4811106313602522
$ws.Cells.Item(15, 12).Value = 0.2300405548868554
$ws.Cells.Item(15, 14).Value = 1.738892757855929
$ws.Cells.Item(15, 15).Value = 3.445598450817442

$ws.Cells.Item(16, 2).Value = 0.8202683763628045
$ws.Cells.Item(16, 3).Value = 0.1462202276487403
$ws.Cells.Item(16, 5).Value = 0.1160661580148563
$ws.Cells.Item(16, 6).Value = 0.6671810134426437
$ws.Cells.Item(16, 7).Value = 0.002462195150467323
$ws.Cells.Item(16, 11).Value = 0.4555628432796937
$ws.Cells.Item(16, 12).Value = 0.2241530405048735
$ws.Cells.Item(16, 14).Value = 1.750451834923183
$ws.Cells.Item(16, 15).Value = 3.457196850497439

$ws.Cells.Item(17, 2).Value = 0.8011643287289303
$ws.Cells.Item(17, 3).Value = 0.1462924940618038
$ws.Cells.Item(17, 5).Value = 0.1155890319815249
$ws.Cells.Item(17, 6).Value = 0.6400460337125793
$ws.Cells.Item(17, 7).Value = 0.002463088310252166
$ws.Cells.Item(17, 11).Value = 0.4398988332804947
$ws.Cells.Item(17, 12).Value = 0.2205684950554598
$ws.Cells.Item(17, 14).Value = 1.757702016404563
$ws.Cells.Item(17, 15).Value = 3.464787732538156

$ws.Cells.Item(18, 2).Value = 0.7901985271647618
$ws.Cells.Item(18, 3).Value = 0.1463358322824213
$ws.Cells.Item(18, 5).Value = 0.1153202520881713
$ws.Cells.Item(18, 6).Value = 0.6244449056556647
$ws.Cells.Item(18, 7).Value = 0.002463609269848444
$ws.Cells.Item(18, 11).Value = 0.4308926827333153
$ws.Cells.Item(18, 12).Value = 0.2185168043484964
$ws.Cells.Item(18, 14).Value = 1.761930577222858
$ws.Cells.Item(18, 15).Value = 3.469328580370444

$ws.Cells.Item(19, 2).Value = 0.7864895535597043
$ws.Cells.Item(19, 3).Value = 0.1463508108646856
$ws.Cells.Item(19, 5).Value = 0.1152302190167624
$ws.Cells.Item(19, 6).Value = 0.619163680173358
$ws.Cells.Item(19, 7).Value = 0.002463786902386392
$ws.Cells.Item(19, 11).Value = 0.4278439500186266
$ws.Cells.Item(19, 12).Value = 0.2178238651433446
$ws.Cells.Item(19, 14).Value = 1.763372333200525
$ws.Cells.Item(19, 15).Value = 3.47089604875822

$ws.Cells.Item(20, 2).Value = 0.8031956800154489
$ws.Cells.Item(20, 3).Value = 0.146284617869064
$ws.Cells.Item(20, 5).Value = 0.11563923824135
$ws.Cells.Item(20, 6).Value = 0.642933953830422
$ws.Cells.Item(20, 7).Value = 0.002462992483236393
$ws.Cells.Item(20, 11).Value = 0.4415659484460832
$ws.Cells.Item(20, 12).Value = 0.2209490370750586
$ws.Cells.Item(20, 14).Value = 1.756924172895333
$ws.Cells.Item(20, 15).Value = 3.463961581296559

$ws.Cells.Item(21, 2).Value = 0.8595913359027634
$ws.Cells.Item(21, 3).Value = 0.1460833533202042
$ws.Cells.Item(21, 5).Value = 0.1170809428882045
$ws.Cells.Item(21, 6).Value = 0.7228739723491628
$ws.Cells.Item(21, 7).Value = 0.002460409703288455
$ws.Cells.Item(21, 11).Value = 0.4877091324263745
$ws.Cells.Item(21, 12).Value = 0.2315689972680417
$ws.Cells.Item(21, 14).Value = 1.735957612723688
$ws.Cells.Item(21, 15).Value = 3.442751039410723

$ws.Cells.Item(22, 2).Value = 0.8966433461434917
$ws.Cells.Item(22, 3).Value = 0.1459674705272285
$ws.Cells.Item(22, 5).Value = 0.1180733334622239
$ws.Cells.Item(22, 6).Value = 0.7751780083420101
$ws.Cells.Item(22, 7).Value = 0.002458785816018447
$ws.Cells.Item(22, 11).Value = 0.5178929015096969
$ws.Cells.Item(22, 12).Value = 0.238598581289736
$ws.Cells.Item(22, 14).Value = 1.722775374907734
$ws.Cells.Item(22, 15).Value = 3.430445542661801

$ws.Cells.Item(23, 2).Value = 0.8768504530505368
$ws.Cells.Item(23, 3).Value = 0.1460278930991059
$ws.Cells.Item(23, 5).Value = 0.1175390948102191
$ws.Cells.Item(23, 6).Value = 0.7472568307830727
$ws.Cells.Item(23, 7).Value = 0.002459646665745792
$ws.Cells.Item(23, 11).Value = 0.5017809983261827
$ws.Cells.Item(23, 12).Value = 0.2348386626950401
$ws.Cells.Item(23, 14).Value = 1.729763338756547
$ws.Cells.Item(23, 15).Value = 3.436870745344493

$ws.Cells.Item(24, 2).Value = 0.8022772519410353
$ws.Cells.Item(24, 3).Value = 0.146288173113998
$ws.Cells.Item(24, 5).Value = 0.1156165227762997
$ws.Cells.Item(24, 6).Value = 0.6416283278902171
$ws.Cells.Item(24, 7).Value = 0.002463035783546622
$ws.Cells.Item(24, 11).Value = 0.4408122478295411
$ws.Cells.Item(24, 12).Value = 0.2207769656391179
$ws.Cells.Item(24, 14).Value = 1.757275648005693
$ws.Cells.Item(24, 15).Value = 3.464334533715117

$ws.Cells.Item(25, 2).Value = 0.7226946406428567
$ws.Cells.Item(25, 3).Value = 0.1466367705971336
$ws.Cells.Item(25, 5).Value = 0.1137588265504377
$ws.Cells.Item(25, 6).Value = 0.5279251897347166
$ws.Cells.Item(25, 7).Value = 0.002466968412801555
$ws.Cells.Item(25, 11).Value = 0.3751783986207329
$ws.Cells.Item(25, 12).Value = 0.2059934913472858
$ws.Cells.Item(25, 14).Value = 1.789186728447358
$ws.Cells.Item(25, 15).Value = 3.500640310129825
